# The document's single-row/single-cell table is replaced with a new
# 2-row x 3-column table carrying the real report header content
# (Chinese labels + mail-merge-style placeholders), and the old
# red-on-dark-shaded "text" cell is gone entirely (no more w:tcPr/w:shd).
#
# The cleanest, least error-prone way to reproduce this exact OOXML
# shape via the Word object model is to hand Word a WordOpenXML package
# for the new table and have it replace the whole document body with
# InsertXML -- this guarantees no leftover tcPr/shading cruft survives
# on the cells (something the individual Shading / Font property setters
# cannot reliably clear), while w:sectPr is left untouched by Word.

$d = $word.ActiveDocument

$tableXml = '<w:tbl>' +
  '<w:tr>' +
    '<w:tc><w:p><w:r><w:rPr><w:color w:val="00FF00"/></w:rPr>' +
      '<w:t>列印人員:&amp;name&amp;(&amp;nowdate$)</w:t></w:r></w:p></w:tc>' +
    '<w:tc><w:p><w:r><w:rPr><w:b/><w:sz w:val="30"/></w:rPr>' +
      '<w:t>國泰敦南健檢中心</w:t></w:r></w:p></w:tc>' +
    '<w:tc><w:p><w:r><w:rPr><w:b/><w:sz w:val="40"/></w:rPr>' +
      '<w:t>*$chartno$*</w:t></w:r></w:p></w:tc>' +
  '</w:tr>' +
  '<w:tr>' +
    '<w:tc><w:p><w:r><w:t/></w:r></w:p></w:tc>' +
    '<w:tc><w:p><w:r><w:rPr><w:b/></w:rPr>' +
      '<w:t>敦南健檢健檢報告(院內)</w:t></w:r></w:p></w:tc>' +
    '<w:tc><w:p><w:r><w:t/></w:r></w:p></w:tc>' +
  '</w:tr>' +
'</w:tbl>'

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' + $tableXml + '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$d.Content.InsertXML($packageXml)
